$d = $word.ActiveDocument

# Merge the split "<id>", "p057v_1"/"p057v_2", "</id>" runs back into single runs,
# as would happen when the <id> value was re-downloaded/retyped as one piece of text.
$d.Content.Find.Execute("<id>p057v_1</id>", $false, $false, $false, $false, $false, $true, 1, $false, "<id>p057v_1</id>", 2)
$d.Content.Find.Execute("<id>p057v_2</id>", $false, $false, $false, $false, $false, $true, 1, $false, "<id>p057v_2</id>", 2)
